# PVRP_Comparison_Results.xlsx update
# Adds new "CPSC_464" (J) / "%_Worse" (K) result columns for rows 3-6, 8-44
# (row 7 already had values which are refreshed with tweaked-algorithm output),
# reflecting an algorithm tweak for some of the more intense PVRP formats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @{Row=3;  J=1057.8831787109375;  K=101.65519990677421},
  @{Row=4;  J=2329.040283203125;   K=76.05964933841761},
  @{Row=5;  J=1036.7662353515625;  K=97.62985805405307},
  @{Row=6;  J=1630.29736328125;    K=95.18441722113474},
  @{Row=7;  J=4118.6328125;        K=103.08940441027815},
  @{Row=8;  J=1630.29736328125;    K=95.14002792282602},
  @{Row=9;  J=1740.581787109375;   K=110.6884773899551},
  @{Row=10; J=4622.36572265625;    K=127.23819397076173},
  @{Row=11; J=1740.581787109375;   K=110.6884773899551},
  @{Row=12; J=3843.853515625;      K=141.22837338008725},
  @{Row=13; J=1651.6834716796875;  K=112.00979021894176},
  @{Row=14; J=2275.407958984375;   K=90.27059228220011},
  @{Row=15; J=10404.5673828125;    K=196.2896720833262},
  @{Row=16; J=1344.745361328125;   K=40.840527998337336},
  @{Row=17; J=2697.02294921875;    K=44.7988268666783},
  @{Row=18; J=4076.252685546875;   K=41.77283964756799},
  @{Row=19; J=2766.99072265625;    K=73.18045518111407},
  @{Row=20; J=5614.75927734375;    K=78.8316450778182},
  @{Row=21; J=8550.53125;          K=76.87070520484698},
  @{Row=22; J=14552.2802734375;    K=73.91639306639459},
  @{Row=23; J=4158.62353515625;    K=91.58778109177834},
  @{Row=24; J=8446.65234375;       K=101.4008832663718},
  @{Row=25; J=12944.6767578125;    K=101.60818286159162},
  @{Row=26; J=6702.60400390625;    K=81.76750402461991},
  @{Row=27; J=6702.60400390625;    K=77.45135893216445},
  @{Row=28; J=6702.60400390625;    K=76.60181496965342},
  @{Row=29; J=47038.265625;        K=114.23877584714884},
  @{Row=30; J=47038.265625;        K=110.88342802665193},
  @{Row=31; J=47038.265625;        K=107.76756747504953},
  @{Row=32; J=179766.265625;       K=141.41281418092385},
  @{Row=33; J=179766.265625;       K=134.8281933254738},
  @{Row=34; J=179766.265625;       K=130.254431020093},
  @{Row=35; J=3269.283447265625;   K=47.99700533565223},
  @{Row=36; J=7132.0732421875;     K=88.97464666151308},
  @{Row=37; J=9942.7802734375;     K=92.12545092292012},
  @{Row=38; J=12157.337890625;     K=105.53646265678543},
  @{Row=39; J=14563.259765625;     K=120.02371623331496},
  @{Row=40; J=17321.375;           K=109.75063210818979},
  @{Row=41; J=8236.38671875;       K=64.85500243688128},
  @{Row=42; J=14926.22265625;      K=113.54260925904995},
  @{Row=43; J=22524.71875;         K=123.56153353713006},
  @{Row=44; J=26379.75;            K=104.10401511529122}
)

foreach ($item in $data) {
  $ws.Cells.Item($item.Row, 10).Value = $item.J
  $ws.Cells.Item($item.Row, 11).Value = $item.K
}

# The newly populated J/K columns (and the now-wider B-H columns, which grew
# because CPSC_464/%_Worse push the "best" column widths out) were re-fit by
# Excel's bestFit column sizing. Re-apply the closest achievable widths.
$ws.Columns.Item(1).ColumnWidth  = 21.833333333333332
$ws.Columns.Item(2).ColumnWidth  = 8.833333333333334
$ws.Columns.Item(3).ColumnWidth  = 8.833333333333334
$ws.Columns.Item(4).ColumnWidth  = 6.5
$ws.Columns.Item(5).ColumnWidth  = 6.5
$ws.Columns.Item(6).ColumnWidth  = 6.5
$ws.Columns.Item(7).ColumnWidth  = 7.666666666666667
$ws.Columns.Item(8).ColumnWidth  = 8.833333333333334
$ws.Columns.Item(10).ColumnWidth = 17.833333333333332
$ws.Columns.Item(11).ColumnWidth = 15.666666666666666
